$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.774.47'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.251.99'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.85'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '297.06'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +7.13%  '
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.96'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0920'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("E13").Value = '  -1.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.06'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +21.96%  '
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.30'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.591.22'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.252.69'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.723.38'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000106'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.17'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.43'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.51'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +16.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '257.74'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +11.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.92'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.69%  '
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.10'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.26'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +5.64%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '175.23'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.15'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0887'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.68'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.07'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.129'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.22'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0376'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.62%  '
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("E41").Value = '  -6.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.07'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.231'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.93%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.46'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.01%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.52'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '107.77'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +7.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.29'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.69'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.89'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.13%  '
